$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A:A").Insert()
$ws.Range("A1").Value = "IN"
$ws.Range("A2").Formula = "=ROW()-1"
$ws.Range("A3:A31").Formula = "=ROW()-1"
[void]$ws.Range("B20").Select()
